$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price cells that are about to be rewritten to remain text so
# values like "1.001" or "49.80" are not coerced to numbers (matches the
# original inlineStr/text cells). D8, D20, D21 are left untouched since
# their Price value does not change. (Multi-area "A1,A2" refs only apply
# to the first area in this COM layer, so issue one call per block.)
$ws.Range("D2:D7").NumberFormat = "@"
$ws.Range("D9:D19").NumberFormat = "@"
$ws.Range("D22:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.649.30'
$ws.Range("E2").Value = '  +4.60%  '
$ws.Range("D3").Value = '1.824.88'
$ws.Range("E3").Value = '  +5.90%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.56%  '
$ws.Range("D5").Value = '336.94'
$ws.Range("E5").Value = '  +1.04%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.42%  '
$ws.Range("D7").Value = '0.3833'
$ws.Range("E7").Value = '  +2.48%  '
$ws.Range("E8").Value = '  +5.41%  '
$ws.Range("D9").Value = '49.80'
$ws.Range("E9").Value = '  +3.52%  '
$ws.Range("D10").Value = '1.238'
$ws.Range("E10").Value = '  +5.37%  '
$ws.Range("D11").Value = '0.07743'
$ws.Range("E11").Value = '  +5.03%  '
$ws.Range("D12").Value = '1.004'
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("D13").Value = '22.35'
$ws.Range("E13").Value = '  +11.21%  '
$ws.Range("D14").Value = '6.635'
$ws.Range("E14").Value = '  +4.39%  '
$ws.Range("D15").Value = '1.827.10'
$ws.Range("E15").Value = '  +5.66%  '
$ws.Range("D16").Value = '7.195'
$ws.Range("E16").Value = '  +1.86%  '
$ws.Range("D17").Value = '0.00001126'
$ws.Range("E17").Value = '  +5.34%  '
$ws.Range("D18").Value = '0.06725'
$ws.Range("E18").Value = '  +1.20%  '
$ws.Range("D19").Value = '87.32'
$ws.Range("E19").Value = '  +5.82%  '
$ws.Range("E20").Value = '  -0.51%  '
$ws.Range("E21").Value = '  +7.12%  '
$ws.Range("D22").Value = '6.540'
$ws.Range("E22").Value = '  +6.91%  '
$ws.Range("D23").Value = '13.17'
$ws.Range("E23").Value = '  +2.98%  '
$ws.Range("D24").Value = '27.590.19'
$ws.Range("E24").Value = '  +4.46%  '
$ws.Range("D25").Value = '2.476'
$ws.Range("E25").Value = '  +0.77%  '
$ws.Range("D26").Value = '2.661'
$ws.Range("E26").Value = '  +12.20%  '
$ws.Range("D27").Value = '22.20'
$ws.Range("E27").Value = '  +14.55%  '
$ws.Range("D28").Value = '1.482'
$ws.Range("E28").Value = '  +7.44%  '
$ws.Range("D29").Value = '152.51'
$ws.Range("E29").Value = '  -0.74%  '
$ws.Range("D30").Value = '2.032.04'
$ws.Range("E30").Value = '  +5.53%  '
$ws.Range("D31").Value = '135.50'
$ws.Range("E31").Value = '  +3.51%  '
$ws.Range("D32").Value = '6.334'
$ws.Range("E32").Value = '  +6.87%  '
$ws.Range("D33").Value = '4.092'
$ws.Range("E33").Value = '  -1.39%  '
$ws.Range("D34").Value = '13.93'
$ws.Range("E34").Value = '  +10.35%  '
$ws.Range("D35").Value = '0.08827'
$ws.Range("E35").Value = '  +3.06%  '
$ws.Range("D36").Value = '1.701'
$ws.Range("E36").Value = '  +0.58%  '
$ws.Range("D37").Value = '5.616'
$ws.Range("E37").Value = '  +4.55%  '
$ws.Range("D38").Value = '0.7003'
$ws.Range("E38").Value = '  +13.86%  '
$ws.Range("D39").Value = '9.100'
$ws.Range("E39").Value = '  +7.20%  '
$ws.Range("D40").Value = '0.06526'
$ws.Range("E40").Value = '  +5.64%  '
$ws.Range("D41").Value = '0.2261'
$ws.Range("E41").Value = '  +4.77%  '
$ws.Range("D42").Value = '0.02404'
$ws.Range("E42").Value = '  +3.18%  '
$ws.Range("D43").Value = '1.299'
$ws.Range("E43").Value = '  +6.45%  '
$ws.Range("D44").Value = '14.77'
$ws.Range("E44").Value = '  +6.26%  '
$ws.Range("D45").Value = '0.6602'
$ws.Range("E45").Value = '  +11.05%  '
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").Value = '  -0.52%  '
$ws.Range("D47").Value = '3.890'
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("D48").Value = '2.182'
$ws.Range("E48").Value = '  +7.32%  '
$ws.Range("D49").Value = '133.23'
$ws.Range("E49").Value = '  +4.71%  '
$ws.Range("D50").Value = '0.07318'
$ws.Range("E50").Value = '  +1.84%  '
$ws.Range("D51").Value = '81.63'
$ws.Range("E51").Value = '  +6.56%  '
